$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 838.8333
$ws.Range("I19").Value = 752.1667
$ws.Range("J19").Value = 882.1667
$ws.Range("K19").Value = 752.1667
$ws.Range("L19").Value = 882.1667
$ws.Range("M19").Value = -577.1667
$ws.Range("N19").Value = -1232.1667
$ws.Range("H116").Value = 6081
$ws.Range("I116").Value = 6081
$ws.Range("K116").Value = 6081
$ws.Range("M116").Value = -2639
$ws.Range("H137").Value = 1910.3889
$ws.Range("I137").Value = 899
$ws.Range("K137").Value = 2697
$ws.Range("M137").Value = -147
$ws.Range("H138").Value = 5260.5356
$ws.Range("I138").Value = 2020.3334
$ws.Range("J138").Value = 5649.36
$ws.Range("K138").Value = 6061.0002
$ws.Range("L138").Value = 16948.08
$ws.Range("M138").Value = -921.0002000000004
$ws.Range("N138").Value = -27228.08

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 46000
$ws.Range("J37").Value = 46000
$ws.Range("L37").Value = 46000
$ws.Range("N37").Value = -46546
$ws.Range("H74").Value = 1879.1428
$ws.Range("I74").Value = 1541
$ws.Range("K74").Value = 1541
$ws.Range("M74").Value = -667
$ws.Range("H77").Value = 1879.1428
$ws.Range("I77").Value = 1541
$ws.Range("K77").Value = 7705
$ws.Range("M77").Value = -3337
$ws.Range("H110").Value = 1496
$ws.Range("I110").Value = 1343.4546
$ws.Range("J110").Value = 2055.3333
$ws.Range("K110").Value = 1343.4546
$ws.Range("L110").Value = 2055.3333
$ws.Range("M110").Value = 701.5454
$ws.Range("N110").Value = -6145.3333
$ws.Range("H132").Value = 595
$ws.Range("I132").Value = 595
$ws.Range("K132").Value = 1785
$ws.Range("M132").Value = 745

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 7356.6665
$ws.Range("I36").Value = 7356.6665
$ws.Range("K36").Value = 7356.6665
$ws.Range("M36").Value = -6822.6665
$ws.Range("H94").Value = 2447.077
$ws.Range("J94").Value = 634.3333
$ws.Range("L94").Value = 634.3333
$ws.Range("N94").Value = -1536.3333
$ws.Range("H107").Value = 889.8570999999999
$ws.Range("I107").Value = 873.1667
$ws.Range("K107").Value = 873.1667
$ws.Range("M107").Value = 1046.8333
$ws.Range("H134").Value = 2654.739
$ws.Range("I134").Value = 2482
$ws.Range("J134").Value = 3276.6
$ws.Range("K134").Value = 7446
$ws.Range("L134").Value = 9829.799999999999
$ws.Range("M134").Value = -4911
$ws.Range("N134").Value = -14899.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 959.5714
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H86").Value = 4580.9
$ws.Range("I86").Value = 5588.3335
$ws.Range("J86").Value = 4149.143
$ws.Range("K86").Value = 5588.3335
$ws.Range("L86").Value = 4149.143
$ws.Range("M86").Value = -4465.3335
$ws.Range("N86").Value = -6395.143
$ws.Range("H89").Value = 4580.9
$ws.Range("I89").Value = 5588.3335
$ws.Range("J89").Value = 4149.143
$ws.Range("K89").Value = 27941.6675
$ws.Range("L89").Value = 20745.715
$ws.Range("M89").Value = -22325.6675
$ws.Range("N89").Value = -31977.715
$ws.Range("H113").Value = 959.5714
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 3913.7273
$ws.Range("I122").Value = 3150
$ws.Range("K122").Value = 9450
$ws.Range("M122").Value = -7000
$ws.Range("H132").Value = 4197.467
$ws.Range("I132").Value = 3896.8
$ws.Range("J132").Value = 4798.8
$ws.Range("K132").Value = 11690.4
$ws.Range("L132").Value = 14396.4
$ws.Range("M132").Value = -9160.400000000001
$ws.Range("N132").Value = -19456.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1459.375
$ws.Range("I68").Value = 1405.5
$ws.Range("J68").Value = 1467.0714
$ws.Range("K68").Value = 4216.5
$ws.Range("L68").Value = 4401.2142
$ws.Range("M68").Value = -3405.5
$ws.Range("N68").Value = -6023.2142
$ws.Range("H71").Value = 1459.375
$ws.Range("I71").Value = 1405.5
$ws.Range("J71").Value = 1467.0714
$ws.Range("K71").Value = 12649.5
$ws.Range("L71").Value = 13203.6426
$ws.Range("M71").Value = -8593.5
$ws.Range("N71").Value = -21315.6426
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("N82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").ClearContents()
$ws.Range("N85").Value = 0
$ws.Range("H103").Value = 5199
$ws.Range("J103").Value = 5399
$ws.Range("L103").Value = 16197
$ws.Range("N103").Value = -17955
$ws.Range("H129").Value = 3850.75
$ws.Range("I129").Value = 1145.625
$ws.Range("K129").Value = 3436.875
$ws.Range("M129").Value = 1563.125
$ws.Range("H139").Value = 6489.2
$ws.Range("I139").Value = 4202.4546
$ws.Range("J139").Value = 12777.75
$ws.Range("K139").Value = 12607.3638
$ws.Range("L139").Value = 38333.25
$ws.Range("M139").Value = -7467.363799999999
$ws.Range("N139").Value = -48613.25
$ws.Range("H140").Value = 2082.7646
$ws.Range("I140").Value = 2082.7646
$ws.Range("K140").Value = 6248.293799999999
$ws.Range("M140").Value = -1068.293799999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 11477273
$ws.Range("I11").Value = 15500000
$ws.Range("K11").Value = 15500000
$ws.Range("M11").Value = -15499861
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H21").Value = 5000
$ws.Range("I21").Value = 5000
$ws.Range("K21").Value = 5000
$ws.Range("M21").Value = -4827
$ws.Range("H30").Value = 5000
$ws.Range("I30").Value = 5000
$ws.Range("K30").Value = 5000
$ws.Range("M30").Value = -4895
$ws.Range("H80").Value = 8457.833000000001
$ws.Range("I80").Value = 4408.2
$ws.Range("K80").Value = 4408.2
$ws.Range("M80").Value = -3410.2
$ws.Range("H83").Value = 8457.833000000001
$ws.Range("I83").Value = 4408.2
$ws.Range("K83").Value = 22041
$ws.Range("M83").Value = -17049
$ws.Range("H107").Value = 1371
$ws.Range("I107").Value = 890
$ws.Range("J107").Value = 2333
$ws.Range("K107").Value = 890
$ws.Range("L107").Value = 2333
$ws.Range("M107").Value = 1030
$ws.Range("N107").Value = -6173
$ws.Range("H122").Value = 4129.6665
$ws.Range("I122").Value = 4129.6665
$ws.Range("K122").Value = 12388.9995
$ws.Range("M122").Value = -9938.999500000002
$ws.Range("H132").Value = 2631.923
$ws.Range("I132").Value = 2057.611
$ws.Range("J132").Value = 3924.125
$ws.Range("K132").Value = 6172.833
$ws.Range("L132").Value = 11772.375
$ws.Range("M132").Value = -3642.833
$ws.Range("N132").Value = -16832.375

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5693.3335
$ws.Range("I22").Value = 5694.5
$ws.Range("J22").Value = 5691
$ws.Range("K22").Value = 5694.5
$ws.Range("L22").Value = 5691
$ws.Range("M22").Value = -5399.5
$ws.Range("N22").Value = -6281
$ws.Range("H23").Value = 59666.332
$ws.Range("I23").Value = 39499.5
$ws.Range("K23").Value = 39499.5
$ws.Range("M23").Value = -39269.5
$ws.Range("H27").Value = 5693.3335
$ws.Range("I27").Value = 5694.5
$ws.Range("J27").Value = 5691
$ws.Range("K27").Value = 5694.5
$ws.Range("L27").Value = 5691
$ws.Range("M27").Value = -5587.5
$ws.Range("N27").Value = -5905
$ws.Range("H34").Value = 27500
$ws.Range("J34").Value = 25000
$ws.Range("L34").Value = 25000
$ws.Range("N34").Value = -25344
$ws.Range("H40").Value = 4831.222
$ws.Range("I40").Value = 3996
$ws.Range("J40").Value = 5499.4
$ws.Range("K40").Value = 3996
$ws.Range("L40").Value = 5499.4
$ws.Range("M40").Value = -3860
$ws.Range("N40").Value = -5771.4
$ws.Range("H46").Value = 3028.2856
$ws.Range("J46").Value = 3333.1667
$ws.Range("L46").Value = 3333.1667
$ws.Range("N46").Value = -3709.1667
$ws.Range("H74").Value = 39975
$ws.Range("J74").Value = 39975
$ws.Range("L74").Value = 39975
$ws.Range("N74").Value = -41971
$ws.Range("H77").Value = 39975
$ws.Range("J77").Value = 39975
$ws.Range("L77").Value = 119925
$ws.Range("N77").Value = -129909
$ws.Range("H93").Value = 2946.0715
$ws.Range("I93").Value = 2386.4546
$ws.Range("K93").Value = 2386.4546
$ws.Range("M93").Value = -1138.4546
$ws.Range("H100").Value = 4642.2856
$ws.Range("I100").Value = 4749.5
$ws.Range("K100").Value = 4749.5
$ws.Range("M100").Value = -4208.5
$ws.Range("H132").Value = 4306.2
$ws.Range("I132").Value = 3177.111
$ws.Range("K132").Value = 9531.332999999999
$ws.Range("M132").Value = -7001.332999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 841.5625
$ws.Range("I113").Value = 781.4
$ws.Range("K113").Value = 2344.2
$ws.Range("M113").Value = -174.1999999999998
$ws.Range("H122").Value = 14582
$ws.Range("I122").Value = 19498
$ws.Range("J122").Value = 4750
$ws.Range("K122").Value = 58494
$ws.Range("L122").Value = 14250
$ws.Range("M122").Value = -56044
$ws.Range("N122").Value = -19150
$ws.Range("H126").Value = 4658.684
$ws.Range("I126").Value = 2895.889
$ws.Range("K126").Value = 8687.667000000001
$ws.Range("M126").Value = -6217.667000000001
$ws.Range("H132").Value = 6838.3706
$ws.Range("I132").Value = 3865.5454
$ws.Range("J132").Value = 19918.8
$ws.Range("K132").Value = 11596.6362
$ws.Range("L132").Value = 59756.39999999999
$ws.Range("M132").Value = -9066.636200000001
$ws.Range("N132").Value = -64816.39999999999
